$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Cells.Item($row, 2).Value = "Beste klant,`nBedankt voor je e-mail. Helaas kan ik je in dit geval niet verder helpen met je bestelling van 200 stuks M8-bouten RVS voor Van Dijk. Ik raad je aan om contact op te nemen met onze verkoopafdeling of een van onze vertegenwoordigers, zodat zij je verder kunnen assisteren met het plaatsen van deze bestelling.`nMocht je nog andere vragen hebben of hulp nodig hebben, laat het ons gerust weten.`nMet vriendelijke groet,`n[Naam] E-mailassistent bij [Bedrijfsnaam]"
$ws.Cells.Item($row, 3).Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Cells.Item($row, 4).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 5).Value = "Bestelling / Levering"
$ws.Cells.Item($row, 6).Value = "2025-07-31 22:01:55"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

$ws.Rows.Item($row).AutoFit()
